# Apply edits to the Data.xlsx test workbook:
#  - ProductDetails sheet: update row 2 (Purse -> pens, 5 -> 2)
#    and remove row 3 (Mobiles / 4 / Newest Arrivals) entirely.
#  - ReferenceData sheet content is unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ProductDetails")

# Update row 2 values
$ws.Range("A2").Value = "pens"
$ws.Range("B2").Value = 2

# Remove row 3 (shifts nothing below it up, it's the last row)
$ws.Rows.Item(3).Delete()

$wb.Save()
